# Apply the component/value corrections described in the commit:
# "Poprawiony schemat oraz wartości komponentów"
# (Fixed schematic and component values)

$wb = $excel.ActiveWorkbook

# --- Sheet "Wartości" ---
$wartosci = $wb.Worksheets.Item("Wartości")

# Unit for the chosen inductance value corrected from uF to uH
$wartosci.Range("D29").Value = [char]0x00B5 + "H"

# KIND (ripple coefficient) limit lowered from 30% to 20%, value 0.3 -> 0.2
$wartosci.Range("C12").Value = 0.2
$wartosci.Range("F12").Value = "współczynnik tętnienia (limit 20% tętnienia)"

# R6 / RLK_COMP value changed from 47 to 22
$wartosci.Range("C20").Value = 22

# --- Sheet "Obliczenia" ---
$obliczenia = $wb.Worksheets.Item("Obliczenia")

# L - chosen inductance value 15 -> 22
$obliczenia.Range("E5").Value = 22

# C7 / CCOUT - chosen output capacitance value 33 -> 44
$obliczenia.Range("E10").Value = 44

# R2 value 18 -> 15
$obliczenia.Range("E15").Value = 15

# Output capacitor description updated to reflect doubled capacitance
$obliczenia.Range("G10").Value = "pojemność kondensatora wyjściowego x2"

$excel.Calculate()

# --- Restore / update selection state on each sheet, matching the author's
#     last interaction while reviewing the corrected schematic ---
$dokumentacja = $wb.Worksheets.Item("Dokumentacja")
$dokumentacja.Activate()
$dokumentacja.Range("B3").Select()

$obliczenia.Activate()
$obliczenia.Range("G14:G15").Select()

$wartosci.Activate()
$wartosci.Range("M21").Select()
